$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2, F2, G2
$ws.Cells.Item(2, 4).Value = 63488
$ws.Cells.Item(2, 5).Value = 1251266030203
$ws.Cells.Item(2, 6).Value = 16951389959
$ws.Cells.Item(2, 7).Value = 0.65895

# Row 3: update D3, E3, F3, G3
$ws.Cells.Item(3, 4).Value = 3286.75
$ws.Cells.Item(3, 5).Value = 401721209804
$ws.Cells.Item(3, 6).Value = 13502862795
$ws.Cells.Item(3, 7).Value = 4.76128

# Row 4: update D4, E4, F4, G4
$ws.Cells.Item(4, 4).Value = 0.999579
$ws.Cells.Item(4, 5).Value = 110542916254
$ws.Cells.Item(4, 6).Value = 18832345482
$ws.Cells.Item(4, 7).Value = 0.02277

# Row 5: update D5, E5, F5, G5
$ws.Cells.Item(5, 4).Value = 601.6900000000001
$ws.Cells.Item(5, 5).Value = 92675069960
$ws.Cells.Item(5, 6).Value = 721687671
$ws.Cells.Item(5, 7).Value = 2.27929

# Row 6: update D6, E6, F6, G6
$ws.Cells.Item(6, 4).Value = 141.38
$ws.Cells.Item(6, 5).Value = 63295584877
$ws.Cells.Item(6, 6).Value = 2380296592
$ws.Cells.Item(6, 7).Value = 3.20119

# Row 7: update D7, E7, F7, G7
$ws.Cells.Item(7, 4).Value = 0.999716
$ws.Cells.Item(7, 5).Value = 33507244349
$ws.Cells.Item(7, 6).Value = 4449551523
$ws.Cells.Item(7, 7).Value = -0.03693

# Row 8: update D8, E8, F8, G8
$ws.Cells.Item(8, 4).Value = 3286.35
$ws.Cells.Item(8, 5).Value = 30765593133
$ws.Cells.Item(8, 6).Value = 92147431
$ws.Cells.Item(8, 7).Value = 4.89726

# Row 9: update D9, E9, F9, G9
$ws.Cells.Item(9, 4).Value = 0.517728
$ws.Cells.Item(9, 5).Value = 28616045471
$ws.Cells.Item(9, 6).Value = 516211197
$ws.Cells.Item(9, 7).Value = 0.40678

# Row 10: update D10, E10, F10, G10
$ws.Cells.Item(10, 4).Value = 0.148723
$ws.Cells.Item(10, 5).Value = 21443343857
$ws.Cells.Item(10, 6).Value = 702088757
$ws.Cells.Item(10, 7).Value = 2.51809

# Row 11: update D11, E11, F11, G11
$ws.Cells.Item(11, 4).Value = 5.45
$ws.Cells.Item(11, 5).Value = 18937900907
$ws.Cells.Item(11, 6).Value = 163757771
$ws.Cells.Item(11, 7).Value = 3.73112

# Row 12: update D12, E12, F12, G12
$ws.Cells.Item(12, 4).Value = 0.467783
$ws.Cells.Item(12, 5).Value = 16535064095
$ws.Cells.Item(12, 6).Value = 261054090
$ws.Cells.Item(12, 7).Value = 2.67533

# Row 13: update D13, E13, F13, G13
$ws.Cells.Item(13, 4).Value = 0.00002456
$ws.Cells.Item(13, 5).Value = 14486163800
$ws.Cells.Item(13, 6).Value = 344333433
$ws.Cells.Item(13, 7).Value = 0.715

# Row 14: update D14, E14, F14, G14
$ws.Cells.Item(14, 4).Value = 34.4
$ws.Cells.Item(14, 5).Value = 13034887033
$ws.Cells.Item(14, 6).Value = 310883405
$ws.Cells.Item(14, 7).Value = 0.95145

# Row 15: update B15, C15, D15, E15, F15, G15
$ws.Cells.Item(15, 2).Value = 'TRX'
$ws.Cells.Item(15, 3).Value = 'TRON'
$ws.Cells.Item(15, 4).Value = 0.121092
$ws.Cells.Item(15, 5).Value = 10614772091
$ws.Cells.Item(15, 6).Value = 251070783
$ws.Cells.Item(15, 7).Value = 0.79

# Row 16: update B16, C16, D16, E16, F16, G16
$ws.Cells.Item(16, 2).Value = 'WBTC'
$ws.Cells.Item(16, 3).Value = 'Wrapped Bitcoin'
$ws.Cells.Item(16, 4).Value = 63545
$ws.Cells.Item(16, 5).Value = 9878528213
$ws.Cells.Item(16, 6).Value = 284359962
$ws.Cells.Item(16, 7).Value = 0.83538

# Row 17: update B17, C17, D17, E17, F17, G17
$ws.Cells.Item(17, 2).Value = 'BCH'
$ws.Cells.Item(17, 3).Value = 'Bitcoin Cash'
$ws.Cells.Item(17, 4).Value = 477.98
$ws.Cells.Item(17, 5).Value = 9424278716
$ws.Cells.Item(17, 6).Value = 231262089
$ws.Cells.Item(17, 7).Value = 1.34032

# Row 18: update D18, E18, F18, G18
$ws.Cells.Item(18, 4).Value = 6.82
$ws.Cells.Item(18, 5).Value = 9258124401
$ws.Cells.Item(18, 6).Value = 140412822
$ws.Cells.Item(18, 7).Value = 2.40286

# Row 19: update D19, E19, F19, G19
$ws.Cells.Item(19, 4).Value = 14.04
$ws.Cells.Item(19, 5).Value = 8254974336
$ws.Cells.Item(19, 6).Value = 281480770
$ws.Cells.Item(19, 7).Value = -0.4885

# Row 20: update B20, C20, D20, E20, F20, G20
$ws.Cells.Item(20, 2).Value = 'NEAR'
$ws.Cells.Item(20, 3).Value = 'NEAR Protocol'
$ws.Cells.Item(20, 4).Value = 7.23
$ws.Cells.Item(20, 5).Value = 7747521248
$ws.Cells.Item(20, 6).Value = 670497251
$ws.Cells.Item(20, 7).Value = 5.45859

# Row 21: update B21, C21, D21, E21, F21, G21
$ws.Cells.Item(21, 2).Value = 'MATIC'
$ws.Cells.Item(21, 3).Value = 'Polygon'
$ws.Cells.Item(21, 4).Value = 0.725945
$ws.Cells.Item(21, 5).Value = 6748832234
$ws.Cells.Item(21, 6).Value = 283000054
$ws.Cells.Item(21, 7).Value = 4.23283

# Row 22: update B22, C22, D22, E22, F22, G22
$ws.Cells.Item(22, 2).Value = 'LTC'
$ws.Cells.Item(22, 3).Value = 'Litecoin'
$ws.Cells.Item(22, 4).Value = 84.01000000000001
$ws.Cells.Item(22, 5).Value = 6260460660
$ws.Cells.Item(22, 6).Value = 352925134
$ws.Cells.Item(22, 7).Value = -1.59217

# Row 23: update B23, C23, D23, E23, F23, G23
$ws.Cells.Item(23, 2).Value = 'ICP'
$ws.Cells.Item(23, 3).Value = 'Internet Computer'
$ws.Cells.Item(23, 4).Value = 13.53
$ws.Cells.Item(23, 5).Value = 6259219232
$ws.Cells.Item(23, 6).Value = 108246132
$ws.Cells.Item(23, 7).Value = 4.32893

# Row 24: update D24, E24, F24, G24
$ws.Cells.Item(24, 4).Value = 8
$ws.Cells.Item(24, 5).Value = 6037159757
$ws.Cells.Item(24, 6).Value = 105914827
$ws.Cells.Item(24, 7).Value = 4.45935

# Row 25: update D25, E25, F25, G25
$ws.Cells.Item(25, 4).Value = 5.77
$ws.Cells.Item(25, 5).Value = 5353560880
$ws.Cells.Item(25, 6).Value = 1115193
$ws.Cells.Item(25, 7).Value = -0.18576

# Row 26: update D26, E26, F26, G26
$ws.Cells.Item(26, 4).Value = 0.9990869999999999
$ws.Cells.Item(26, 5).Value = 5300051791
$ws.Cells.Item(26, 6).Value = 341533742
$ws.Cells.Item(26, 7).Value = -0.04717

# Row 27: update B27, C27, D27, E27, F27, G27
$ws.Cells.Item(27, 2).Value = 'FDUSD'
$ws.Cells.Item(27, 3).Value = 'First Digital USD'
$ws.Cells.Item(27, 4).Value = 0.999356
$ws.Cells.Item(27, 5).Value = 4421963328
$ws.Cells.Item(27, 6).Value = 4223605304
$ws.Cells.Item(27, 7).Value = -0.41116

# Row 28: update B28, C28, D28, E28, F28, G28
$ws.Cells.Item(28, 2).Value = 'ETC'
$ws.Cells.Item(28, 3).Value = 'Ethereum Classic'
$ws.Cells.Item(28, 4).Value = 28.64
$ws.Cells.Item(28, 5).Value = 4208131568
$ws.Cells.Item(28, 6).Value = 292584727
$ws.Cells.Item(28, 7).Value = 7.47523

# Row 29: update B29, C29, D29, E29, F29, G29
$ws.Cells.Item(29, 2).Value = 'APT'
$ws.Cells.Item(29, 3).Value = 'Aptos'
$ws.Cells.Item(29, 4).Value = 8.99
$ws.Cells.Item(29, 5).Value = 3845874138
$ws.Cells.Item(29, 6).Value = 94904639
$ws.Cells.Item(29, 7).Value = 2.51416

# Row 30: update B30, C30, D30, E30, F30, G30
$ws.Cells.Item(30, 2).Value = 'HBAR'
$ws.Cells.Item(30, 3).Value = 'Hedera'
$ws.Cells.Item(30, 4).Value = 0.104187
$ws.Cells.Item(30, 5).Value = 3730540254
$ws.Cells.Item(30, 6).Value = 127272948
$ws.Cells.Item(30, 7).Value = -2.67316

# Row 31: update B31, C31, D31, E31, F31, G31
$ws.Cells.Item(31, 2).Value = 'STX'
$ws.Cells.Item(31, 3).Value = 'Stacks'
$ws.Cells.Item(31, 4).Value = 2.52
$ws.Cells.Item(31, 5).Value = 3670743372
$ws.Cells.Item(31, 6).Value = 50466822
$ws.Cells.Item(31, 7).Value = -0.30715

# Row 32: update B32, C32, D32, E32, F32, G32
$ws.Cells.Item(32, 2).Value = 'MNT'
$ws.Cells.Item(32, 3).Value = 'Mantle'
$ws.Cells.Item(32, 4).Value = 1.1
$ws.Cells.Item(32, 5).Value = 3593450469
$ws.Cells.Item(32, 6).Value = 38972817
$ws.Cells.Item(32, 7).Value = 3.00505

# Row 33: update B33, C33, D33, E33, F33, G33
$ws.Cells.Item(33, 2).Value = 'CRO'
$ws.Cells.Item(33, 3).Value = 'Cronos'
$ws.Cells.Item(33, 4).Value = 0.128041
$ws.Cells.Item(33, 5).Value = 3418998297
$ws.Cells.Item(33, 6).Value = 12572539
$ws.Cells.Item(33, 7).Value = 3.33036

# Row 34: update B34, C34, D34, E34, F34, G34
$ws.Cells.Item(34, 2).Value = 'XLM'
$ws.Cells.Item(34, 3).Value = 'Stellar'
$ws.Cells.Item(34, 4).Value = 0.114372
$ws.Cells.Item(34, 5).Value = 3313398133
$ws.Cells.Item(34, 6).Value = 49343280
$ws.Cells.Item(34, 7).Value = 1.43716

# Row 35: update D35, E35, F35, G35
$ws.Cells.Item(35, 4).Value = 8.31
$ws.Cells.Item(35, 5).Value = 3249455817
$ws.Cells.Item(35, 6).Value = 86327410
$ws.Cells.Item(35, 7).Value = 1.30525

# Row 36: update B36, C36, D36, E36, F36, G36
$ws.Cells.Item(36, 2).Value = 'FIL'
$ws.Cells.Item(36, 3).Value = 'Filecoin'
$ws.Cells.Item(36, 4).Value = 5.95
$ws.Cells.Item(36, 5).Value = 3243569658
$ws.Cells.Item(36, 6).Value = 139021506
$ws.Cells.Item(36, 7).Value = 3.2917

# Row 37: update B37, C37, D37, E37, F37, G37
$ws.Cells.Item(37, 2).Value = 'EZETH'
$ws.Cells.Item(37, 3).Value = 'Renzo Restaked ETH'
$ws.Cells.Item(37, 4).Value = 3237.11
$ws.Cells.Item(37, 5).Value = 3227342167
$ws.Cells.Item(37, 6).Value = 76829365
$ws.Cells.Item(37, 7).Value = 5.24272

# Row 38: update B38, C38, D38, E38, F38, G38
$ws.Cells.Item(38, 2).Value = 'OKB'
$ws.Cells.Item(38, 3).Value = 'OKB'
$ws.Cells.Item(38, 4).Value = 53.13
$ws.Cells.Item(38, 5).Value = 3194753925
$ws.Cells.Item(38, 6).Value = 4785082
$ws.Cells.Item(38, 7).Value = 1.50823

# Row 39: update B39, C39, D39, E39, F39, G39
$ws.Cells.Item(39, 2).Value = 'IMX'
$ws.Cells.Item(39, 3).Value = 'Immutable'
$ws.Cells.Item(39, 4).Value = 2.16
$ws.Cells.Item(39, 5).Value = 3147743787
$ws.Cells.Item(39, 6).Value = 45642754
$ws.Cells.Item(39, 7).Value = 3.17299

# Row 40: update B40, C40, D40, E40, F40, G40
$ws.Cells.Item(40, 2).Value = 'XT'
$ws.Cells.Item(40, 3).Value = 'XT.com'
$ws.Cells.Item(40, 4).Value = 3.13
$ws.Cells.Item(40, 5).Value = 3138034436
$ws.Cells.Item(40, 6).Value = 742336
$ws.Cells.Item(40, 7).Value = 4.62931

# Row 41: update B41, C41, D41, E41, F41, G41
$ws.Cells.Item(41, 2).Value = 'RNDR'
$ws.Cells.Item(41, 3).Value = 'Render'
$ws.Cells.Item(41, 4).Value = 8.09
$ws.Cells.Item(41, 5).Value = 3133954511
$ws.Cells.Item(41, 6).Value = 97679653
$ws.Cells.Item(41, 7).Value = 2.16761

# Row 42: update B42, C42, D42, E42, F42, G42
$ws.Cells.Item(42, 2).Value = 'PEPE'
$ws.Cells.Item(42, 3).Value = 'Pepe'
$ws.Cells.Item(42, 4).Value = 0.00000737
$ws.Cells.Item(42, 5).Value = 3098542896
$ws.Cells.Item(42, 6).Value = 669852942
$ws.Cells.Item(42, 7).Value = 6.81877

# Row 43: update B43, C43, D43, E43, F43, G43
$ws.Cells.Item(43, 2).Value = 'ARB'
$ws.Cells.Item(43, 3).Value = 'Arbitrum'
$ws.Cells.Item(43, 4).Value = 1.12
$ws.Cells.Item(43, 5).Value = 2970633387
$ws.Cells.Item(43, 6).Value = 348154838
$ws.Cells.Item(43, 7).Value = 6.60358

# Row 44: update B44, C44, D44, E44, F44, G44
$ws.Cells.Item(44, 2).Value = 'VET'
$ws.Cells.Item(44, 3).Value = 'VeChain'
$ws.Cells.Item(44, 4).Value = 0.03964638
$ws.Cells.Item(44, 5).Value = 2888436480
$ws.Cells.Item(44, 6).Value = 48118811
$ws.Cells.Item(44, 7).Value = 2.84993

# Row 45: update B45, C45, D45, E45, F45, G45
$ws.Cells.Item(45, 2).Value = 'TAO'
$ws.Cells.Item(45, 3).Value = 'Bittensor'
$ws.Cells.Item(45, 4).Value = 427.75
$ws.Cells.Item(45, 5).Value = 2852288918
$ws.Cells.Item(45, 6).Value = 19784582
$ws.Cells.Item(45, 7).Value = 2.94756

# Row 46: update B46, C46, D46, E46, F46, G46
$ws.Cells.Item(46, 2).Value = 'MKR'
$ws.Cells.Item(46, 3).Value = 'Maker'
$ws.Cells.Item(46, 4).Value = 3061.68
$ws.Cells.Item(46, 5).Value = 2832605504
$ws.Cells.Item(46, 6).Value = 95132555
$ws.Cells.Item(46, 7).Value = 5.36623

# Row 47: update B47, C47, D47, E47, F47, G47
$ws.Cells.Item(47, 2).Value = 'OP'
$ws.Cells.Item(47, 3).Value = 'Optimism'
$ws.Cells.Item(47, 4).Value = 2.64
$ws.Cells.Item(47, 5).Value = 2776157258
$ws.Cells.Item(47, 6).Value = 431748611
$ws.Cells.Item(47, 7).Value = 14.87225

# Row 48: update B48, C48, D48, E48, F48, G48
$ws.Cells.Item(48, 2).Value = 'WIF'
$ws.Cells.Item(48, 3).Value = 'dogwifhat'
$ws.Cells.Item(48, 4).Value = 2.73
$ws.Cells.Item(48, 5).Value = 2723643176
$ws.Cells.Item(48, 6).Value = 317126823
$ws.Cells.Item(48, 7).Value = 1.04446

# Row 49: update B49, C49, D49, E49, F49, G49
$ws.Cells.Item(49, 2).Value = 'WEETH'
$ws.Cells.Item(49, 3).Value = 'Wrapped eETH'
$ws.Cells.Item(49, 4).Value = 3405.45
$ws.Cells.Item(49, 5).Value = 2716518222
$ws.Cells.Item(49, 6).Value = 33378922
$ws.Cells.Item(49, 7).Value = 4.93554

# Row 50: update B50, C50, D50, E50, F50, G50
$ws.Cells.Item(50, 2).Value = 'KAS'
$ws.Cells.Item(50, 3).Value = 'Kaspa'
$ws.Cells.Item(50, 4).Value = 0.111848
$ws.Cells.Item(50, 5).Value = 2584743947
$ws.Cells.Item(50, 6).Value = 44428073
$ws.Cells.Item(50, 7).Value = 0.15864

# Row 51: update B51, C51, D51, E51, F51, G51
$ws.Cells.Item(51, 2).Value = 'GRT'
$ws.Cells.Item(51, 3).Value = 'The Graph'
$ws.Cells.Item(51, 4).Value = 0.262865
$ws.Cells.Item(51, 5).Value = 2496616408
$ws.Cells.Item(51, 6).Value = 67095824
$ws.Cells.Item(51, 7).Value = 0.88644
